
function Find-ParaIndex($text) {
    $d = $word.ActiveDocument
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -eq ($text + [char]13)) {
            return $i
        }
    }
    return -1
}

function Replace-ParaText($oldText, $newText) {
    $idx = Find-ParaIndex $oldText
    if ($idx -eq -1) {
        Write-Output ("REPLACE ANCHOR NOT FOUND: " + $oldText)
        return
    }
    $d = $word.ActiveDocument
    $p = $d.Paragraphs.Item($idx)
    $p.Range.Text = $newText
}

function Insert-ParaBefore($anchorText, $newText) {
    $idx = Find-ParaIndex $anchorText
    if ($idx -eq -1) {
        Write-Output ("INSERT-BEFORE ANCHOR NOT FOUND: " + $anchorText)
        return
    }
    $d = $word.ActiveDocument
    $p = $d.Paragraphs.Item($idx)
    $p.Range.InsertParagraphBefore()
    $newPara = $d.Paragraphs.Item($idx)
    $newPara.Range.Text = $newText
}

function Insert-ParaAfter($anchorText, $newText) {
    $idx = Find-ParaIndex $anchorText
    if ($idx -eq -1) {
        Write-Output ("INSERT-AFTER ANCHOR NOT FOUND: " + $anchorText)
        return
    }
    $d = $word.ActiveDocument
    $p = $d.Paragraphs.Item($idx)
    $p.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($idx + 1)
    $newPara.Range.Text = $newText
}

function Set-StandardRunFormatting($paraText) {
    # Some anchor paragraphs end in a stray empty run that only carries
    # <w:rtl/>, so a paragraph inserted adjacent to it can inherit
    # incomplete run formatting. Force the standard Calibri/12pt (24
    # half-point) formatting used throughout this document onto the
    # paragraph's run so the OOXML run properties come out complete.
    $idx = Find-ParaIndex $paraText
    if ($idx -eq -1) {
        Write-Output ("SET-FONT ANCHOR NOT FOUND: " + $paraText)
        return
    }
    $d = $word.ActiveDocument
    $rng = $d.Paragraphs.Item($idx).Range
    $f = $rng.Font
    $f.Name = "Calibri"
    $f.NameFarEast = "Calibri"
    $f.NameOther = "Calibri"
    $f.NameAscii = "Calibri"
    $f.NameBi = "Calibri"
    $f.Size = 12
    $f.SizeBi = 12
}

# --- Step A: text replacements ---
Replace-ParaText "Asher (neutral expressionless): What took you so long?" "Asher (neutral neutral): What took you so long?"
Replace-ParaText "Asher (neutral grinning): Like Prim?" "Asher (neutral smirk): Like Prim?"
Replace-ParaText "Asher (excited grinning): Everyone knows how you two are going steady." "Asher (neutral smirk): Everyone knows how you two are going steady."
Replace-ParaText "Asher (neutral grinning): Yup." "Asher (neutral smiling_eyes_closed): Yup."
Replace-ParaText "Asher (neutral playful): I dunno…" "Asher (neutral thinking): I dunno…"
Replace-ParaText "He holds out for a few more seconds before breaking out in laughter, attracting the attention of all of our classmates. Between fits of chuckles he hands be his phone, and upon closer inspection I realize that a certain hyperactive first year recently sent out a few irresponsibly false messages." "He holds out for a few more seconds before breaking out in laughter, attracting the attention of all of our classmates. Between fits of chuckles he hands me his phone, and upon closer inspection I realize that a certain hyperactive first year recently sent out a few irresponsibly false messages."
Replace-ParaText "A few of Asher’s friends take interest in the conversation and join in to interrogate me, which is a little uncomfortable since I’m not used to all the attention. Fortunately, they eventually lose interest and leave me be." "A few of Asher’s friends take interest in the conversation and surround me, all of them curious as to what’s gotten him so worked up. Which is a little uncomfortable since I’m not used to all the attention, but fortunately they eventually lose interest and leave me be."
Replace-ParaText "Asher (laughing recovering): Sorry, sorry. It was too funny to resist though." "Asher (neutral smiling_eyes_closed): Sorry, sorry. It was too funny to resist though."
Replace-ParaText "Asher (neutral smiling): And you tutor her." "Asher (neutral smirk): And you tutor her."
Replace-ParaText "Asher (neutral thinking): Doesn’t it? How else would you measure if you’ve gotten closer to someone?" "Asher (neutral curious): Doesn’t it, though? How else would you measure if you’ve gotten closer to someone?"
Replace-ParaText "I think about it for a second, and after a few moments I realize that Asher’s as right as ever." "I think about it for a second, and after a few moments I realize that I can’t come up with a counterargument. As usual."
Replace-ParaText "Asher (neutral grinning): You sure?" "Asher (neutral smirk): You sure?"
Replace-ParaText "Asher starts laughing again, causing me to lean back into my chair and laugh." "Asher starts laughing again, causing me to lean back into my chair and sigh."

# --- Step B: paragraph insertions ---
Insert-ParaBefore "He lets out a chuckle." "Asher (neutral hehe):"
Insert-ParaBefore "He holds out for a few more seconds before breaking out in laughter, attracting the attention of all of our classmates. Between fits of chuckles he hands me his phone, and upon closer inspection I realize that a certain hyperactive first year recently sent out a few irresponsibly false messages." "Asher (laughing laughing):"
Insert-ParaAfter "A few of Asher’s friends take interest in the conversation and surround me, all of them curious as to what’s gotten him so worked up. Which is a little uncomfortable since I’m not used to all the attention, but fortunately they eventually lose interest and leave me be." "Asher (laughing recovering):"
Insert-ParaBefore "Pro: But it’s still not like that." "Asher (neutral smiling):"
Insert-ParaBefore "Asher starts laughing again, causing me to lean back into my chair and sigh." "Asher (laughing laughing):"
Insert-ParaAfter "Asher starts laughing again, causing me to lean back into my chair and sigh." "Asher (exit):"
Insert-ParaAfter "Asher (exit):" "Mick (arms_crossed annoyed):"
Insert-ParaBefore "Was that…" "Mick (exit):"
Replace-ParaText "Asher continues to poke fun at me until lunch ends regardless. But despite that, in the back of my mind I can’t help but wonder why Mick would peer into our classroom like that." "Asher continues to poke fun at me until lunch ends regardless. But despite that, in the back of my mind I can’t help but wonder why Mick would peer into our classroom like that with such a hostile look. I could be reading too much into it, but…"
Insert-ParaBefore "Asher continues to poke fun at me until lunch ends regardless. But despite that, in the back of my mind I can’t help but wonder why Mick would peer into our classroom like that with such a hostile look. I could be reading too much into it, but…" "Asher (exit):"
Insert-ParaAfter "Asher continues to poke fun at me until lunch ends regardless. But despite that, in the back of my mind I can’t help but wonder why Mick would peer into our classroom like that with such a hostile look. I could be reading too much into it, but…" "…was he looking at me?"
Set-StandardRunFormatting "Asher (exit):"
Set-StandardRunFormatting "…was he looking at me?"
